# Update cryptocurrency price (D) and 1-hour volume change (E) columns
# with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.596.31"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.472.29"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0862"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "2.850.98"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").Value = "2.464.22"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("D18").Value = "41.547.89"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0763"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  +4.64%  "
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.64%  "
$ws.Range("D43").Value = "1.984.29"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("E47").Value = "  +5.73%  "
$ws.Range("D48").Value = "2.708.34"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.83%  "
